$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 12

# Start the new row with the same formatting (center aligned) as the
# previous data row, so every cell we touch already carries that style.
$ws.Range("A11:C11").Copy() | Out-Null
$ws.Range("A$row`:C$row").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$cellA = $ws.Cells.Item($row, 1)
$cellB = $ws.Cells.Item($row, 2)
$cellC = $ws.Cells.Item($row, 3)

# Date column: force text so "2025/11/21" is stored as a literal string
# instead of being auto-converted into a date serial number.
$cellA.NumberFormat = "@"
$cellA.Value = "2025/11/21"

# Game column (plain text)
$cellB.Value = "逃离鸭科夫"

# ModCount column (numeric)
$cellC.Value = 1211

# Re-apply the reference row's formatting so the date cell ends up with the
# exact same style as its neighbours (center aligned, default number format)
# rather than keeping the temporary text format used above.
$ws.Range("A11:C11").Copy() | Out-Null
$ws.Range("A$row`:C$row").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
